# addressbook_contact_overview.xlsx - "adb xls export improvements" edit
#
# - the salutation placeholder moves from the first data column to use the
#   keyField() twig helper (translates the raw salutation value through the
#   Addressbook "contactSalutation" keyfield definition) instead of the raw
#   record value
# - the now-unused column I (left over from an earlier layout) is cleared
# - header / data row heights are tightened up
# - the stray manual column break is removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- salutation cell: route the value through keyField() -------------------
$ws.Range("A6").Value = '${ROW}${twig:keyField(''Addressbook'',''contactSalutation'',record.salutation)}'

# --- drop the no-longer-used column I (was spans 1:9, now 1:8) -------------
$ws.Range("I1:I6").Clear()

# --- tighten row heights for the header/data rows --------------------------
$ws.Rows(5).RowHeight = 23
$ws.Rows(6).RowHeight = 20

# --- remove the manual column page break ------------------------------------
$ws.ResetAllPageBreaks()

# --- rename the default cell style typo "Stand." -> "Standard" -------------
foreach ($s in $wb.Styles) {
  if ($s.Name -eq "Stand.") {
    $s.Name = "Standard"
  }
}

# --- reset selection back to the top-left cell ------------------------------
$ws.Range("A1").Select() | Out-Null
